# Apply the new Orders rows (22-27) and refresh the Summary totals string,
# mirroring the addition of new flower line items to the order sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New order rows to append. Column A (PackageID) and F (Number) hold
# digit-only strings in this workbook, so force Text format on those
# cells before assigning the value - otherwise Excel would coerce them
# to numeric cells instead of text-as-number cells.
$newRows = @(
    @{ Row = 22; A = $null; C = "590_洋牡丹粉_undefined_undefined_1bunch"; F = "20" },
    @{ Row = 23; A = $null; C = "585_洋牡丹红_undefined_undefined_1bunch"; F = "10" },
    @{ Row = 24; A = "34";  C = "419_松虫草红_scabiosa watermelon_undefined_1bunch"; F = "50" },
    @{ Row = 25; A = $null; C = "512_松虫草粉_scabiosa pink_undefined_1bunch"; F = "50" },
    @{ Row = 26; A = $null; C = "514_松虫草紫_scabiosa purple_undefined_1bunch"; F = "20" },
    @{ Row = 27; A = $null; C = "447_黄金球_craspedia_undefined_1bunch"; F = "20" }
)

foreach ($r in $newRows) {
    if ($r.A -ne $null) {
        $cellA = $ws.Cells.Item($r.Row, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $r.A
    }

    $ws.Cells.Item($r.Row, 3).Value = $r.C

    $cellF = $ws.Cells.Item($r.Row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $r.F
}

# The Summary sheet keeps a running concatenation of every Number value
# from the Orders sheet in G2; extend it with the six new quantities.
$g2 = $summary.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "015205205803030205151510155532201920201050502020"
